# Add a "location_name" column (E) and link the full roster of teams with
# their founding year, trainer and logo, one row per club.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Formatting first -------------------------------------------------
# New header cell E1 should look like the other header cells (A1:D1).
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)

# New data cells (E2:E7) and the new rows (4:7, columns A:D) should look
# like the existing data rows (A2:D3).
$ws.Range("A2").Copy()
$ws.Range("E2:E7").PasteSpecial(-4122)
$ws.Range("A2:D2").Copy()
$ws.Range("A4:D7").PasteSpecial(-4122)

# Column D (team_logo) needs to be a bit wider now that it also has to fit
# the new logo file names.
$ws.Columns("D").ColumnWidth = 16.42

# --- Values -------------------------------------------------------------
# Write every cell in row-major order (row by row, column A -> E) so the
# workbook's shared-string table is rebuilt in the same order it would be
# if the rows were entered top-to-bottom, left-to-right.

# Row 1 - header
$ws.Cells.Item(1,1).Value = "team_name"
$ws.Cells.Item(1,2).Value = "team_foundingDate"
$ws.Cells.Item(1,3).Value = "team_trainer"
$ws.Cells.Item(1,4).Value = "team_logo"
$ws.Cells.Item(1,5).Value = "location_name"

# Row 2 - Hells Teddies
$ws.Cells.Item(2,1).Value = "Hells Teddies"
$ws.Cells.Item(2,2).Value = 2013
$ws.Cells.Item(2,3).Value = "test trainer"
$ws.Cells.Item(2,4).Value = "/Team-logos/logo1.png"
$ws.Cells.Item(2,5).Value = "Ulaanbaatar"

# Row 3 - PPL
$ws.Cells.Item(3,1).Value = "PPL"
$ws.Cells.Item(3,2).Value = 2018
$ws.Cells.Item(3,3).Value = "prfkps"
$ws.Cells.Item(3,4).Value = "/Team-logos/logo3.png"
$ws.Cells.Item(3,5).Value = "Lisma"

# Row 4 - Elements
$ws.Cells.Item(4,1).Value = "Elements"
$ws.Cells.Item(4,2).Value = 2002
$ws.Cells.Item(4,3).Value = "florian"
$ws.Cells.Item(4,4).Value = "/Team-logos/logo2.jpg"
$ws.Cells.Item(4,5).Value = "Gusswerk"

# Row 5 - Valantic
$ws.Cells.Item(5,1).Value = "Valantic"
$ws.Cells.Item(5,2).Value = 2018
$ws.Cells.Item(5,3).Value = "johanna"
$ws.Cells.Item(5,4).Value = "/Team-logos/logo3.png"
$ws.Cells.Item(5,5).Value = "Salzburg"

# Row 6 - New Team
$ws.Cells.Item(6,1).Value = "New Team"
$ws.Cells.Item(6,2).Value = 2024
$ws.Cells.Item(6,3).Value = "new trainer"
$ws.Cells.Item(6,4).Value = "/Team-logos/logo4.jpg"
$ws.Cells.Item(6,5).Value = "Glasgow"

# Row 7 - Alchimiste
$ws.Cells.Item(7,1).Value = "Alchimiste"
$ws.Cells.Item(7,2).Value = 2009
$ws.Cells.Item(7,3).Value = "dirk"
$ws.Cells.Item(7,4).Value = "/Team-logos/logo5.jpg"
$ws.Cells.Item(7,5).Value = "Bergstrasse"
